# Rubric tracker update: mark the "Pending" rubric rows (16-29) as
# "Completed", fill in their completion dates, and record the "Completed
# as Specified." result note for each (presentation slides + updated
# rubric were added and these items got wrapped up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> completion date (Excel serial date number)
$rowDates = [ordered]@{
    16 = 44020
    17 = 44020
    18 = 44020
    19 = 44021
    20 = 44021
    21 = 44021
    22 = 44021
    23 = 44020
    24 = 44021
    25 = 44021
    26 = 44021
    27 = 44021
    28 = 44021
    29 = 44021
}

foreach ($r in $rowDates.Keys) {
    $ws.Range("D$r").Value = "Completed"
    $ws.Range("E$r").Value = $rowDates[$r]
    $ws.Range("F$r").Value = "Completed as Specified."
}

# Match Excel's re-wrap of these two rows once their Result text is filled in.
$ws.Rows.Item(17).RowHeight = 86.4
$ws.Rows.Item(24).RowHeight = 86.4

# Leave the view scrolled back to the top of the sheet, as after the edit.
$ws.Range("E1").Select()
